$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.726580142974854
$ws.Range("B1").Value = 2.31868052482605
$ws.Range("C1").Value = 2.399679183959961
$ws.Range("D1").Value = 2.695592880249023
$ws.Range("E1").Value = 3.423704862594604
